# Auto-generated PowerShell-style Excel COM-interop script
# Applies "İş Takip Güncellemesi" date shifts (-1 day) plus a couple of
# status text updates, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("İş Takip Listesi")
$ws2 = $wb.Worksheets.Item("Güncelleme")

# ---- Sheet: İş Takip Listesi (columns J = İŞE BAŞLAMA/YER TESLİMİ, K = İHALE BİTİŞ TARİHİ) ----
$sheet1Updates = @(
    @{Row=2; J="2025-08-22"; K="2026-01-22"},
    @{Row=3; J="2025-08-22"; K="2026-01-22"},
    @{Row=4; J="2025-08-22"; K="2026-01-22"},
    @{Row=5; J="2025-08-22"; K="2026-01-22"},
    @{Row=6; J="2025-08-22"; K="2026-01-22"},
    @{Row=7; J="2025-08-22"; K="2026-01-22"},
    @{Row=8; J="2025-08-22"; K="2026-01-22"},
    @{Row=9; J="2025-08-22"; K="2026-01-22"},
    @{Row=10; J="2025-08-22"; K="2026-01-22"},
    @{Row=33; J="2025-08-24"; K="2026-01-24"},
    @{Row=34; J="2025-08-24"; K="2026-01-24"},
    @{Row=35; J="2025-08-24"; K="2026-01-24"},
    @{Row=36; J="2025-08-24"; K="2026-01-24"},
    @{Row=37; J="2025-08-24"; K="2026-01-24"},
    @{Row=38; J="2025-08-24"; K="2026-01-24"},
    @{Row=39; J="2025-08-24"; K="2026-01-24"},
    @{Row=40; J="2025-08-24"; K="2026-01-24"},
    @{Row=41; J="2025-08-24"; K="2026-01-24"},
    @{Row=42; J="2025-08-24"; K="2026-01-24"},
    @{Row=43; J="2025-08-24"; K="2026-01-24"},
    @{Row=44; J="2025-08-24"; K="2026-01-24"},
    @{Row=45; J="2025-08-24"; K="2026-01-24"},
    @{Row=46; J="2025-08-24"; K="2026-01-24"},
    @{Row=47; J="2025-08-24"; K="2026-01-24"},
    @{Row=48; J="2025-08-24"; K="2026-01-24"},
    @{Row=49; J="2025-08-24"; K="2026-01-24"},
    @{Row=50; J="2025-08-24"; K="2026-01-24"},
    @{Row=51; J="2025-08-24"; K="2026-01-24"},
    @{Row=52; J="2025-08-24"; K="2026-01-24"},
    @{Row=53; J="2025-08-24"; K="2026-01-24"},
    @{Row=54; J="2025-08-24"; K="2026-01-24"},
    @{Row=55; J="2025-08-24"; K="2026-01-24"},
    @{Row=56; J="2025-08-24"; K="2026-01-24"},
    @{Row=57; J="2025-08-24"; K="2026-01-24"},
    @{Row=58; J="2025-08-24"; K="2026-01-24"},
    @{Row=59; J="2025-08-24"; K="2026-01-24"},
    @{Row=60; J="2025-08-24"; K="2026-01-24"},
    @{Row=61; J="2025-08-24"; K="2026-01-24"},
    @{Row=62; J="2025-08-24"; K="2026-01-24"},
    @{Row=63; J="2025-08-24"; K="2026-01-24"},
    @{Row=64; J="2025-08-24"; K="2026-01-24"},
    @{Row=65; J="2025-08-24"; K="2026-01-24"},
    @{Row=66; J="2025-08-24"; K="2026-01-24"},
    @{Row=67; J="2025-08-24"; K="2026-01-24"},
    @{Row=68; J="2025-08-24"; K="2026-01-24"},
    @{Row=69; J="2025-08-24"; K="2026-01-24"},
    @{Row=70; J="2025-08-24"; K="2026-01-24"},
    @{Row=71; J="2025-08-24"; K="2026-01-24"},
    @{Row=72; J="2025-08-24"; K="2026-01-24"},
    @{Row=73; J="2025-08-24"; K="2026-01-24"},
    @{Row=74; J="2025-08-24"; K="2026-01-24"},
    @{Row=75; J="2025-08-24"; K="2026-01-24"},
    @{Row=76; J="2025-08-24"; K="2026-01-24"},
    @{Row=77; J="2025-08-24"; K="2026-01-24"},
    @{Row=78; J="2025-08-24"; K="2026-01-24"},
    @{Row=79; J="2025-08-24"; K="2026-01-24"},
    @{Row=80; J="2025-08-24"; K="2026-01-24"},
    @{Row=81; J="2025-08-24"; K="2026-01-24"},
    @{Row=82; J="2025-08-24"; K="2026-01-24"},
    @{Row=83; J="2025-08-24"; K="2026-01-24"},
    @{Row=84; J="2025-08-24"; K="2026-01-24"},
    @{Row=85; J="2025-08-24"; K="2026-01-24"},
    @{Row=86; J="2025-08-24"; K="2026-01-24"},
    @{Row=87; J="2025-08-24"; K="2026-01-24"},
    @{Row=88; J="2025-08-24"; K="2026-01-24"},
    @{Row=89; J="2025-08-24"; K="2026-01-24"},
    @{Row=90; J="2025-08-24"; K="2026-01-24"},
    @{Row=91; J="2025-08-24"; K="2026-01-24"},
    @{Row=92; J="2025-08-24"; K="2026-01-24"},
    @{Row=93; J="2025-08-24"; K="2026-01-24"},
    @{Row=94; J="2025-08-24"; K="2026-01-24"},
    @{Row=95; J="2024-06-22"; K="2025-08-16"},
    @{Row=96; J="2024-06-22"; K="2025-08-16"},
    @{Row=97; J="2024-06-22"; K="2025-08-16"},
    @{Row=98; J="2024-06-22"; K="2025-08-16"},
    @{Row=99; J="2024-06-22"; K="2025-08-16"},
    @{Row=100; J="2024-06-22"; K="2025-08-16"},
    @{Row=101; J="2024-06-22"; K="2025-08-16"},
    @{Row=102; J="2024-06-22"; K="2025-08-16"},
    @{Row=103; J="2024-06-22"; K="2025-08-16"},
    @{Row=104; J="2024-06-22"; K="2025-08-16"},
    @{Row=105; J="2024-06-22"; K="2025-08-16"},
    @{Row=106; J="2024-06-22"; K="2025-08-16"},
    @{Row=107; J="2024-06-22"; K="2025-08-16"},
    @{Row=108; J="2024-06-22"; K="2025-08-16"},
    @{Row=109; J="2024-06-22"; K="2025-08-16"},
    @{Row=110; J="2024-06-22"; K="2025-08-16"},
    @{Row=111; J="2024-06-22"; K="2025-08-16"},
    @{Row=112; J="2024-06-22"; K="2025-08-16"},
    @{Row=113; J="2024-06-22"; K="2025-08-16"},
    @{Row=114; J="2024-06-22"; K="2025-08-16"},
    @{Row=115; J="2024-06-22"; K="2025-08-16"},
    @{Row=116; J="2024-06-22"; K="2025-08-16"},
    @{Row=117; J="2024-06-22"; K="2025-08-16"},
    @{Row=118; J="2024-06-22"; K="2025-08-16"},
    @{Row=119; J="2024-06-22"; K="2025-08-16"},
    @{Row=120; J="2024-06-22"; K="2025-08-16"},
    @{Row=121; J="2024-06-22"; K="2025-08-16"},
    @{Row=122; J="2024-06-22"; K="2025-08-16"}
)

foreach ($u in $sheet1Updates) {
    $ws1.Cells.Item($u.Row, 10).NumberFormat = "@"
    $ws1.Cells.Item($u.Row, 11).NumberFormat = "@"
    $ws1.Cells.Item($u.Row, 10).Value = $u.J
    $ws1.Cells.Item($u.Row, 11).Value = $u.K
}

# ---- Sheet: Güncelleme (various date columns I, J, N, P) ----
$sheet2Updates = @(
    @{Row=2; Col=10; Value="2024-09-28"},
    @{Row=2; Col=14; Value="2025-06-01"},
    @{Row=2; Col=16; Value="2025-08-18"},
    @{Row=3; Col=10; Value="2024-12-30"},
    @{Row=3; Col=14; Value="2025-09-19"},
    @{Row=4; Col=10; Value="2024-11-03"},
    @{Row=4; Col=14; Value="2025-04-26"},
    @{Row=4; Col=16; Value="2025-07-21"},
    @{Row=5; Col=9; Value="2025-04-29"},
    @{Row=6; Col=10; Value="2025-12-10"},
    @{Row=6; Col=14; Value="2025-08-29"},
    @{Row=7; Col=9; Value="2024-12-30"},
    @{Row=7; Col=10; Value="2024-12-30"},
    @{Row=8; Col=10; Value="2024-12-18"},
    @{Row=8; Col=14; Value="2025-05-19"},
    @{Row=8; Col=16; Value="2025-06-21"},
    @{Row=9; Col=9; Value="2025-08-15"},
    @{Row=9; Col=10; Value="2025-01-31"},
    @{Row=10; Col=10; Value="2024-11-29"},
    @{Row=10; Col=14; Value="2025-09-09"},
    @{Row=11; Col=9; Value="2025-06-06"},
    @{Row=11; Col=10; Value="2025-01-12"},
    @{Row=11; Col=14; Value="2025-09-29"},
    @{Row=12; Col=10; Value="2024-12-10"},
    @{Row=12; Col=14; Value="2025-08-19"},
    @{Row=13; Col=10; Value="2025-02-07"},
    @{Row=14; Col=10; Value="2025-12-06"},
    @{Row=15; Col=10; Value="2025-02-26"},
    @{Row=15; Col=14; Value="2025-09-16"},
    @{Row=16; Col=10; Value="2024-10-25"},
    @{Row=16; Col=14; Value="2025-04-04"},
    @{Row=16; Col=16; Value="2025-06-21"},
    @{Row=17; Col=10; Value="2024-11-10"},
    @{Row=18; Col=10; Value="2025-04-19"},
    @{Row=19; Col=9; Value="2025-06-07"},
    @{Row=19; Col=10; Value="2025-02-26"},
    @{Row=19; Col=14; Value="2025-10-06"},
    @{Row=20; Col=10; Value="2025-02-07"},
    @{Row=21; Col=10; Value="2024-12-01"},
    @{Row=22; Col=10; Value="2024-12-01"},
    @{Row=23; Col=10; Value="2025-02-08"},
    @{Row=24; Col=9; Value="2025-08-05"},
    @{Row=25; Col=10; Value="2025-01-03"},
    @{Row=27; Col=10; Value="2025-03-26"},
    @{Row=28; Col=10; Value="2025-01-22"},
    @{Row=29; Col=9; Value="2025-04-13"},
    @{Row=29; Col=10; Value="2025-02-08"}
)

foreach ($u in $sheet2Updates) {
    $ws2.Cells.Item($u.Row, $u.Col).NumberFormat = "@"
    $ws2.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

# ---- Sheet: Güncelleme - UÇUŞ KIYMETLENDİRME DURUM (column O) flips from blank to "Yapılmadı" ----
$oRows = @(11, 15, 19)
foreach ($r in $oRows) {
    $ws2.Cells.Item($r, 15).Value = "Yapılmadı"
}

